$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 5173.5293
$ws.Range("I100").Value = 5245.8335
$ws.Range("K100").Value = 5245.8335
$ws.Range("M100").Value = -4704.8335
$ws.Range("H125").Value = 999
$ws.Range("I125").Value = 999
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 8991
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -6531
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 9011692
$ws.Range("I122").Value = 12348177
$ws.Range("J122").Value = 3179.9
$ws.Range("K122").Value = 37044531
$ws.Range("L122").Value = 9539.700000000001
$ws.Range("M122").Value = -37042081
$ws.Range("N122").Value = -14439.7
$ws.Range("H132").Value = 32260590
$ws.Range("J132").Value = 3031.4285
$ws.Range("L132").Value = 9094.2855
$ws.Range("N132").Value = -14154.2855
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 49722.9
$ws.Range("I82").Value = 17752.334
$ws.Range("J82").Value = 63424.57
$ws.Range("K82").Value = 17752.334
$ws.Range("L82").Value = 63424.57
$ws.Range("M82").Value = -17369.334
$ws.Range("N82").Value = -64190.57
$ws.Range("H85").Value = 49722.9
$ws.Range("I85").Value = 17752.334
$ws.Range("J85").Value = 63424.57
$ws.Range("K85").Value = 17752.334
$ws.Range("L85").Value = 63424.57
$ws.Range("M85").Value = -16426.334
$ws.Range("N85").Value = -66076.57000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 283.875
$ws.Range("I22").Value = 253
$ws.Range("K22").Value = 253
$ws.Range("M22").Value = 97
$ws.Range("H31").Value = 3331.7925
$ws.Range("I31").Value = 2644.7273
$ws.Range("J31").Value = 3511.738
$ws.Range("K31").Value = 2644.7273
$ws.Range("L31").Value = 3511.738
$ws.Range("M31").Value = -2349.7273
$ws.Range("N31").Value = -4101.737999999999
$ws.Range("H34").Value = 3331.7925
$ws.Range("I34").Value = 2644.7273
$ws.Range("J34").Value = 3511.738
$ws.Range("K34").Value = 2644.7273
$ws.Range("L34").Value = 3511.738
$ws.Range("M34").Value = -2442.7273
$ws.Range("N34").Value = -3915.738
$ws.Range("H99").Value = 2883.05
$ws.Range("I99").Value = 2295.889
$ws.Range("J99").Value = 3363.4546
$ws.Range("K99").Value = 2295.889
$ws.Range("L99").Value = 3363.4546
$ws.Range("M99").Value = -797.8890000000001
$ws.Range("N99").Value = -6359.4546
$ws.Range("H126").Value = 2883.05
$ws.Range("I126").Value = 2295.889
$ws.Range("J126").Value = 3363.4546
$ws.Range("K126").Value = 6887.667
$ws.Range("L126").Value = 10090.3638
$ws.Range("M126").Value = -4417.667
$ws.Range("N126").Value = -15030.3638
$ws.Range("H132").Value = 2293.5454
$ws.Range("I132").Value = 2127.375
$ws.Range("J132").Value = 2736.6667
$ws.Range("K132").Value = 6382.125
$ws.Range("L132").Value = 8210.000100000001
$ws.Range("M132").Value = -3852.125
$ws.Range("N132").Value = -13270.0001
$ws.Range("H134").Value = 1889.1
$ws.Range("I134").Value = 1254.2174
$ws.Range("J134").Value = 3975.1428
$ws.Range("K134").Value = 3762.6522
$ws.Range("L134").Value = 11925.4284
$ws.Range("M134").Value = -1227.6522
$ws.Range("N134").Value = -16995.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 117.5
$ws.Range("I29").Value = 160
$ws.Range("J29").Value = 32.5
$ws.Range("K29").Value = 480
$ws.Range("L29").Value = 97.5
$ws.Range("M29").Value = -203
$ws.Range("N29").Value = -651.5
$ws.Range("H122").Value = 851
$ws.Range("I122").Value = 702
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 6318
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -3868
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 2266.5557
$ws.Range("I132").Value = 1666.6666
$ws.Range("J132").Value = 2566.5
$ws.Range("K132").Value = 14999.9994
$ws.Range("L132").Value = 23098.5
$ws.Range("M132").Value = -12469.9994
$ws.Range("N132").Value = -28158.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 16666.666
$ws.Range("J44").Value = 16666.666
$ws.Range("L44").Value = 16666.666
$ws.Range("N44").Value = -17858.666
$ws.Range("H80").Value = 2633.7896
$ws.Range("I80").Value = 2697.2727
$ws.Range("J80").Value = 2546.5
$ws.Range("K80").Value = 2697.2727
$ws.Range("L80").Value = 2546.5
$ws.Range("M80").Value = -1699.2727
$ws.Range("N80").Value = -4542.5
$ws.Range("H83").Value = 2633.7896
$ws.Range("I83").Value = 2697.2727
$ws.Range("J83").Value = 2546.5
$ws.Range("K83").Value = 13486.3635
$ws.Range("L83").Value = 12732.5
$ws.Range("M83").Value = -8494.363499999999
$ws.Range("N83").Value = -22716.5
$ws.Range("H97").Value = 492.45456
$ws.Range("I97").Value = 492.45456
$ws.Range("K97").Value = 492.45456
$ws.Range("M97").Value = 3.545439999999985
$ws.Range("H132").Value = 4583.8
$ws.Range("I132").Value = 4535
$ws.Range("K132").Value = 13605
$ws.Range("M132").Value = -11075

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 751.1111
$ws.Range("I55").Value = 738.6667
$ws.Range("J55").Value = 776
$ws.Range("K55").Value = 738.6667
$ws.Range("L55").Value = 776
$ws.Range("M55").Value = -565.6667
$ws.Range("N55").Value = -1122
$ws.Range("H61").Value = 12980.25
$ws.Range("I61").Value = 7312.8887
$ws.Range("K61").Value = 7312.8887
$ws.Range("M61").Value = -7110.8887
$ws.Range("H113").Value = 12980.25
$ws.Range("I113").Value = 7312.8887
$ws.Range("K113").Value = 7312.8887
$ws.Range("M113").Value = -5142.8887
$ws.Range("H122").Value = 4272.727
$ws.Range("J122").Value = 5333.3335
$ws.Range("L122").Value = 16000.0005
$ws.Range("N122").Value = -20900.0005
$ws.Range("H136").Value = 1754.04
$ws.Range("I136").Value = 1584.3182
$ws.Range("J136").Value = 2998.6667
$ws.Range("K136").Value = 4752.9546
$ws.Range("L136").Value = 8996.000100000001
$ws.Range("M136").Value = -2202.9546
$ws.Range("N136").Value = -14096.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9096519
$ws.Range("I81").Value = 1955.9166
$ws.Range("K81").Value = 3911.8332
$ws.Range("M81").Value = -2850.8332
$ws.Range("H84").Value = 9096519
$ws.Range("I84").Value = 1955.9166
$ws.Range("K84").Value = 19559.166
$ws.Range("M84").Value = -14255.166
$ws.Range("H100").Value = 3003.5264
$ws.Range("I100").Value = 3971.7856
$ws.Range("K100").Value = 7943.5712
$ws.Range("M100").Value = -7402.5712
$ws.Range("H122").Value = 3599.75
$ws.Range("I122").Value = 3833
$ws.Range("J122").Value = 2900
$ws.Range("K122").Value = 11499
$ws.Range("L122").Value = 8700
$ws.Range("M122").Value = -9049
$ws.Range("N122").Value = -13600
$ws.Range("H132").Value = 2915.889
$ws.Range("I132").Value = 2266.6047
$ws.Range("J132").Value = 5454
$ws.Range("K132").Value = 6799.8141
$ws.Range("L132").Value = 16362
$ws.Range("M132").Value = -4269.8141
$ws.Range("N132").Value = -21422
